$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$ws.Range("C20").Value = "• Mandated nation-wide quarantine for at least 15 days since March 15th, which has been further extended until April 11th. `n• Ban on all non-essential work from March 30 to April 9, inclusive."
$ws.Range("I20").Value = "• Increased sick pay for coronavirus infected workers or those quarantined, from 60 to 75 percent of the regulatory base, paid by the social security budget.`n• Changes in the temporary employment adjustment schemes (ERTEs) in order to avoid outright dismissal by promoting temporary unemployment (the suspension of employment) or reductions in working time. In particular: (i) the temporary employment adjustment schemes have been significantly simplified, (ii) access conditions to these schemes have been extended, (iii) unemployment benefits received under the temporary employment adjustment scheme do not count in terms of consumption of unemployment benefit rights, and (iv) as further detailed under firm-specific measures, employers will be exempt from social contributions during the period of application of the temporary employment adjustment scheme.`n• An extraordinary allowance is provided for self-employed workers affected by the suspension of economic activity (70% of the SSC base, at least for one month).`n• Additional budgetary funds of EUR 300 million to ensure the provision of assistance to dependent persons. `n• Supplemental credit of EUR 25 million to cover meal allowances to ensure the basic access to food for vulnerable children affected by the suspension of educational activity in schools.`n• Additional flexibility for local authorities to use their 2019 budgetary surplus to fund social services and primary assistance to dependent persons (EUR 300 million) `n• The social benefit for energy provision (‘bono social’) will be automatically extended until September 15.`n• Dismissals for reasons realted to Covid19 will not be considered justified (from March 27 till the end of the health crisis)`n`nOther employment and social measures:`n• One-month credit postponement on mortgage payments for the most vulnerable.`n• Broadened scope for protected families in the supply of water and energy, which is ensured for vulnerable groups. Telecommunication services are also guaranteed.`n• Increased flexibility of working conditions, encouragement of telework and adjustment of working times. Reduced working times are permitted for workers having to take care of children, elderly or dependent persons.`n"
$ws.Range("F30").Value = "The government has advised that indoor mass gatherings of 100 people or more and outdoor mass gatherings of more than 500 people should be cancelled. All State-run cultural institutions have closed, as have pubs. On 24 March, the government announced further containment measures, including i) all theatres, clubs, gyms/leisure centres, hairdressers, betting shops, marts, markets, casinos, bingo halls, libraries and other similar outlets to shut; ii) all hotels to limit occupancy to essential non-social and non-tourist reasons; iii) all non-essential retail outlets are to close to members of the public and all other retail outlets are to implement physical distancing`nall cafes and restaurants are to limit supply to take away food or delivery; iv) all sporting events are cancelled, including those behind closed doors`nall playgrounds and holiday/caravan parks will close; v) all places of worship are to restrict numbers entering at any one time to ensure adequate physical distancing`nall organised social indoor and outdoor events of any size are not to take place. On 27 March, the government advised that everybody should stay at home until 12 April 2020, except for the following situations:`ni) to travel to and from work, or for purposes of work, only where the work is an essential health, social care or other essential service and cannot be done from home`nto shop for essential food, beverage and household goods or collect a meal`nii) to attend medical appointments and collect medicines and other health products`nfor vital family reasons, such as providing care to children, elderly or vulnerable people`niii) to take brief individual physical exercise within 2 kilometres of your home.`niv) always observe 2 metres physical distancing`nv) you may take children from your own household outdoors for physical exercise`nvi) farming purposes, that is food production or care of animals."
$ws.Range("G30").Value = "The Irish health system was already experiencing significant challenges at the onset of the COVID-19 outbreak, but a number of measures are being put in place to support the delivery of health services at this critical time. These include agreement with GPs to ensure that cost is not a barrier to accessing services, the utilisation of private hospital capacity, and an expansion in both acute and community capacity. `nThe most recent data available at the onset of the crisis suggested Ireland had 2.8 acute care hospital beds per 1,000 inhabitants. This is low compared with many other European countries. For instance, Germany had six acute care beds per 1,000 inhabitants (the highest across the European countries).   `nThe Government has allocated extra funding of the order of €1 billion (0.5% of GNI*) to the Health Service Executive for improving capacity in the health service, including using private hospital capacity, increasing staffing and overtime, developing home testing and remote management solutions for mildly ill patients and pushing community awareness and preventive actions. The package will also be used to promote greater use of technology and telephone support. Further health-related measures of a similar financial scale are under active consideration.`nTo support public health policy, the Minister of Finance has increased the limit for contactless payments (to EUR 50 per transaction). `n"
$ws.Range("I30").Value = "A package of sick pay and illness/supplementary benefits has been introduced to ensure that employees and the self-employed can abide by medical advice to self-isolate where appropriate, while having their income better protected than under current social welfare rules. For example, the government has waived the requirement for six waiting days for sick pay in respect of medically certified cases of self-isolation in accordance with public health guidelines and has removed the means test requirement for Supplementary Welfare Allowance in respect of medically certified cases of self-isolation. `nThere have also been measures announced to try and encourage businesses to keep workers in employment and to support those who have lost their jobs. On 24 March, the government announced: `ni) a temporary wage subsidy of 70% of take home pay up to a maximum weekly tax free amount of EUR 410 per week to help affected companies keep paying their employees; `nii) workers who have lost their jobs due to the crisis will receive an enhanced emergency COVID-19 Pandemic Unemployment Payment of EUR 350 per week (an increase from EUR 203 per week); `niii) the COVID-19 Illness Payment will be increased to EUR 350 per week; `niv) the self-employed will be eligible for the COVID-19 Pandemic Unemployment Payment of EUR 350 directly from the Department of Employment Affairs and Social Protection.`nThe government is also introducing legislation to prevent both the termination of residential tenancies and any rent increases for the duration of the COVID-19 crisis. Furthermore, the Commission for Regulation of Utilities has issued a moratorium on disconnections of domestic customers for non-payment to the gas and electricity suppliers.`n"
$ws.Range("J30").Value = "The government has announced financial supports for businesses, including; `ni) a €200 million Strategic Banking Corporation of Ireland Working Capital scheme; ii) a €200 million Rescue and Restructuring Scheme available through Enterprise Ireland for vulnerable but viable firms; iii) the Microfinance Ireland COVID-19 Loan includes an option for a moratorium on interest and repayments for the first six months; iv) Local Enterprise Offices in every county will be providing vouchers from €2,500 up to €10,000; v) a Finance in Focus grant of €7,200 will be available to Enterprise Ireland and Údarás na Gaeltachta clients.`nThe government has also agreed with local authorities that they should defer business rates payments due from the most immediately affected businesses, primarily in the retail, hospitality, leisure and childcare sectors, until the end of May. `nThere have also been a variety of taxation measures to alleviate short-term difficulties. For example; i) interest on late payments has been suspended for January/February VAT and both February and March PAYE (Employers) liabilities; ii) all debt enforcement activity is suspended until further notice; iii) current tax clearance status will remain in place for all businesses over the coming months; iv) the Relevant Contract Tax rate review scheduled to take place this month (March) is suspended; v) critical pharmaceutical products and medicines will be given a Customs ‘green routing’ to facilitate uninterrupted importation and supply.`nBanking Initiatives `nFollowing consultation with the Government and the Central Bank of Ireland, the Banking & Payments Federation of Ireland (BPFI) members that includes the five retail Irish banks and main non-banks, have agreed a coordinated approach to assisting customers affected by COVID-19. These measures include:`n• a payment break for up to three months (to be followed by on-going reviews depending on the scale and extent of the situation) for all personal (including personal and buy-to-let mortgages), SME and corporate customers;   `n• seeking to ensure that such a payment break will not adversely impact the borrower’s credit record (the details of this to be discussed with the Central Bank of Ireland); and`n• banks and non-banks deferring court proceedings or any other enforcement action for three months.`n"
$ws.Range("C32").Value = "- As of 9 March all Israelis returning from abroad are ordered into a 14-day quarantine. This measure will be in effect for two weeks initially.`n'- partial lockdown (since 20 March, further tightened 25 March): People are ordered  to stay home except for essential trips (e.g.work, food, medical appointments) (to be reassessed every week). Workplace restrictions with only around 30% of employees being physically  allowed to go to work (since 20 March for 4 weeks). "
$ws.Range("H32").Value = "A first emergency package of NIS 10 billion (0.7% of GDP) announced on 14 March, including NIS 1 billion for health sector, _x000D_NIS 1 billion for special needs (police, fire service, food security),_x000D_ fund for state guaranteed loans of NIS 8 billion for SMEs (off balance sheet measure)_x000D_.`n_x000D_`nFurther measures announced on March 16 include broadened eligibility to unemployment benefits, tax deferals, grant for the self-employed._x000D_`nOverall support measures of around NIS 35 billion (2.5% of GDP).`n."
$ws.Range("I32").Value = "- All employees on involuntary unpaid leave for at least 30 days will be able to claim unemployment benefits after 5 day. Unemployment benefits can be claimed after employment for the least 6 months, from 12 months previously_x000D_. Unemployment benefits for tourist guides.`n'- Grant of up to NIS 6000 for March and April for people over 67 years who quit a job (27 March)`n'- Ministry of Finance instructed banks to allow for a 4-months postponement of mortgage payments (15 March)`n"
$ws.Range("J32").Value = "- Fund for state guaranteed loans of NIS 8 billion for SMEs (increased from previously NIS 4 billion announced on 10 March)_x000D_. Since 29 March, no interest payments for borrower in first year (state pays instead).  `n'- VAT, municipal taxes and social security contribution have been deferred for SMEs until the end April`n'-Grants up to  NIS 6,000 for all small businesses and self employed. Criteria to be determined"
$ws.Range("K32").Value = "On 15 March BOI announced_x000D_`n'- start of purchasing government bonds of various types and maturities in the secondary market (on 23 March specified to totallying NIS 15 billion, 3.5% GDP)_x000D_`n- offering repo transactions to financial institutions with government bonds as collateral (undisclosed amount)_x000D_`n- a list of regulatory leniencies for banks (e.g. higher LTV for all purpose HH loans backed by real estats), aimed at preventing a tightening in credit conditions and to support the economy_x000D_`n-18 March_x000D_ BOI announced to offer banks one-week US dollar swaps up to NIS 15 billion. `n-On 29 March the banking supervisor reduced the regulatory capital requirement for banks by 1 percentatge point"
$ws.Range("C33").Value = "- 23 February: quarantine of municipalities with initial clusters of cases. `n- 25 February: Restrictions introduced to several northern regions. `n- 8 March: Quarantine extended across Lombardy and 14 other northern provinces. `n- 9 March: Nation-wide quarantine. “I stay at home decree”, aims at limiting the movement and contact of people. People are required to stay at home but can leave in case of necessity or for prescribed motives. Reconversion of facilities (e.g., army barracks) to host people in quarantine.`n- 25 March: Movements restrictions reinforced, with fewer exceptions and a limited range of industrial and commercial activities permitted to continue operating, effective 29 March. "
$ws.Range("F33").Value = "- Bars and restaurants along with many other retail trade activities (e.g. shopping centres; indoor and outdoor markets) closed from March 10 until 3 April, and all sporting competitions suspended over the same period along with other public gatherings.`n- All but prescribed essential production activities suspended from March 23, with the list of permitted activities further limited from March 26.`n- On March 30, closures extended from April 3 to 30 April for sports, bars and similar activities. "

# Row 32 height adjustment (Excel will typically autofit on reflow,
# but set explicitly to match the target state)
$ws.Rows.Item(32).RowHeight = 216.75

# Bump "Date Updated" column (B) from 2020-03-29 to 2020-03-30 for all data rows (5-56)
for ($r = 5; $r -le 56; $r++) {
    $ws.Cells.Item($r, 2).Value = 43920
}

